$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.123.73'
$ws.Range('E2').Value = '  -0.71%  '

$ws.Range('D3').Value = '2.648.63'
$ws.Range('E3').Value = '  +0.38%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = "'582.40"
$ws.Range('E5').Value = '  -0.08%  '

$ws.Range('D6').Value = "'156.79"
$ws.Range('E6').Value = '  -0.31%  '

$ws.Range('D7').Value = "'0.628"
$ws.Range('E7').Value = '  -2.57%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').Value = '2.645.17'
$ws.Range('E9').Value = '  +0.38%  '

$ws.Range('E10').Value = '  -3.16%  '

$ws.Range('D11').Value = "'5.83"
$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('E12').Value = '  -1.45%  '

$ws.Range('E13').Value = '  +0.88%  '

$ws.Range('D14').Value = "'28.68"
$ws.Range('E14').Value = '  -0.18%  '

$ws.Range('D15').Value = '3.124.93'

$ws.Range('D16').Value = "'0.0000186"
$ws.Range('E16').Value = '  -0.81%  '

$ws.Range('D17').Value = '63.988.71'
$ws.Range('E17').Value = '  -0.61%  '

$ws.Range('D18').Value = '2.652.05'
$ws.Range('E18').Value = '  +0.68%  '

$ws.Range('D19').Value = "'12.24"
$ws.Range('E19').Value = '  -0.23%  '

$ws.Range('D20').Value = "'7.77"
$ws.Range('E20').Value = '  +3.90%  '

$ws.Range('E21').Value = '  -3.20%  '

$ws.Range('D22').Value = "'346.96"
$ws.Range('E22').Value = '  -0.33%  '

$ws.Range('E23').Value = '  +0.19%  '

$ws.Range('D24').Value = "'68.30"
$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('E25').Value = '  +5.15%  '

$ws.Range('D26').Value = "'0.0000113"
$ws.Range('E26').Value = '  +1.15%  '

$ws.Range('D27').Value = "'9.34"
$ws.Range('E27').Value = '  -1.19%  '

$ws.Range('D28').Value = "'588.77"
$ws.Range('E28').Value = '  -1.14%  '

$ws.Range('E29').Value = '  +1.72%  '

$ws.Range('D30').Value = "'8.27"
$ws.Range('E30').Value = '  +3.01%  '

$ws.Range('E31').Value = '  +0.11%  '

$ws.Range('E32').Value = '  +0.02%  '

$ws.Range('D33').Value = "'2.08"
$ws.Range('E33').Value = '  -0.58%  '

$ws.Range('D34').Value = "'1.76"
$ws.Range('E34').Value = '  +0.33%  '

$ws.Range('D35').Value = "'6.67"
$ws.Range('E35').Value = '  -0.31%  '

$ws.Range('E36').Value = '  +3.36%  '

$ws.Range('E37').Value = '  -2.43%  '

$ws.Range('D38').Value = "'19.83"
$ws.Range('E38').Value = '  -1.30%  '

$ws.Range('E39').Value = '  +0.01%  '

$ws.Range('E40').Value = '  -0.40%  '

$ws.Range('D41').Value = "'151.73"
$ws.Range('E41').Value = '  -1.25%  '

$ws.Range('D42').Value = "'2.59"
$ws.Range('E42').Value = '  +6.78%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'164.09"
$ws.Range('E44').Value = '  +3.32%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'24.55"
$ws.Range('E45').Value = '  +4.64%  '

$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = "'3.93"
$ws.Range('E46').Value = '  -2.41%  '

$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = "'0.0593"
$ws.Range('E47').Value = '  -2.08%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.637"
$ws.Range('E48').Value = '  -0.04%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = "'0.101"
$ws.Range('E49').Value = '  -1.82%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = "'0.0250"
$ws.Range('E50').Value = '  -2.35%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'19.17"
$ws.Range('E51').Value = '  -0.51%  '
